# Analysis dashboard update: add "Date and Time" header row and "Cycle_count" row,
# reorder a few labels, and refresh the computed metric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a brand-new row at the very top for the ride's date/time window.
#    This shifts every existing row down by one (old row 1 -> row 2, etc.)
#    and carries the [hh]:mm:ss style on the "Total time taken" value along
#    with it automatically.
$ws.Rows("1:1").Insert()
$ws.Range("A1").Value = "Date and Time"
$ws.Range("B1").Value = "2024-03-11 13:06:03.902000 to 2024-03-11 13:56:16.218000"

# 2) Remove the old "Maximum BMS Temperature in C" row (now at row 32 after
#    the shift above) - it is dropped entirely in the new layout.
$ws.Rows("32:32").Delete()

# 3) Insert a new row for the battery's cycle count, right before the
#    "Idling time percentage" row (now at row 35).
$ws.Rows("35:35").Insert()
$ws.Range("A35").Value = "Cycle Count of battery"
$ws.Range("B35").Value = 38

# 4) Relabel / reorder the rows whose headings moved or were renamed.
$ws.Range("A7").Value = "Starting SoC (%)"
$ws.Range("B7").Value = 95
$ws.Range("A8").Value = "Ending SoC (%)"
$ws.Range("B8").Value = 25

$ws.Range("A9").Value = "Total distance covered (km)"
$ws.Range("A10").Value = "Total energy consumption(WH/KM)"
$ws.Range("A11").Value = "Total SOC consumed(%)"

$ws.Range("A13").Value = "Peak Power(kW)"
$ws.Range("A14").Value = "Average Power(kW)"
$ws.Range("A15").Value = "Total Energy Regenerated(kWh)"

$ws.Range("A16").Value = "Regenerative Effectiveness(%)"
$ws.Range("B16").Value = 0.04429052093278955

$ws.Range("A17").Value = "Highest Cell Voltage(V)"
$ws.Range("B17").Value = 3.345
$ws.Range("A18").Value = "Lowest Cell Voltage(V)"
$ws.Range("B18").Value = 3

$ws.Range("A19").Value = "Difference in Cell Voltage(V)"
$ws.Range("A20").Value = "Minimum Temperature(C)"
$ws.Range("A21").Value = "Maximum Temperature(C)"

$ws.Range("A22").Value = "Difference in Temperature(C)"
$ws.Range("B22").Value = 10

$ws.Range("A23").Value = "Maximum Fet Temperature-BMS(C)"
$ws.Range("A24").Value = "Maximum Afe Temperature-BMS(C)"
$ws.Range("A25").Value = "Maximum PCB Temperature-BMS(C)"
$ws.Range("A26").Value = "Maximum MCU Temperature(C)"
$ws.Range("A27").Value = "Maximum Motor Temperature(C)"
$ws.Range("A28").Value = "Abnormal Motor Temperature Detected(C)"

$ws.Range("A29").Value = "highest cell temp(C)"
$ws.Range("B29").Value = 47
$ws.Range("A30").Value = "lowest cell temp(C)"
$ws.Range("B30").Value = 37

$ws.Range("A31").Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"

$ws.Range("A32").Value = "Battery Voltage(V)"
$ws.Range("B32").Value = 53

$ws.Range("A33").Value = "Total energy charged(kWh)"
$ws.Range("B33").Value = 1.433210114722222

$ws.Range("A34").Value = "Electricity consumption units(kW)"
$ws.Range("B34").Value = 0.0000001321320679575748

# 5) Refresh the ride-analytics values that were recomputed.
$ws.Range("B36").Value = 23.58145099887604
$ws.Range("A37").Value = "Time spent in 0-10 km/h"
$ws.Range("B37").Value = 11.89949602987564
$ws.Range("B39").Value = 6.8235379427867
$ws.Range("B40").Value = 10.66313766723469
$ws.Range("B41").Value = 11.02933178637468
$ws.Range("B42").Value = 13.40778071861064
$ws.Range("B43").Value = 15.34389616040028

# 6) Append the two new speed-bucket rows at the bottom.
$ws.Range("A44").Value = "Time spent in 70-80 km/h"
$ws.Range("B44").Value = 1.812842173960335
$ws.Range("A45").Value = "Time spent in 80-90 km/h"
$ws.Range("B45").Value = 0
